$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 'osm_id                                                 365436
route                                                     bus
to                                      Amsterdam, Louweshoek
name        Bus 18: Amsterdam Centraal Station => Amsterda...
ref                                                        18
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.900265 52.3797248, 4.89945...
Name: 2, dtype: object'
$ws.Range("A7").Value = 'osm_id                                                 364660
route                                                     bus
to                                    Amsterdam, Nolensstraat
name        Bus 21: Amsterdam Centraal Station => Amsterda...
ref                                                        21
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.900265 52.3797248, 4.89945...
Name: 5, dtype: object'
$ws.Range("A10").Value = 'osm_id                                                8456140
route                                                     bus
to                                     Amsterdam, Noorderpark
name        Bus 34: Amsterdam Olof Palmeplein => Amsterdam...
ref                                                        34
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9401433 52.396837, 4.94000...
Name: 8, dtype: object'
$ws.Range("A11").Value = 'osm_id                                                8456179
route                                                     bus
to                                 Amsterdam, Olof Palmeplein
name        Bus 34: Amsterdam Noorderpark => Amsterdam Olo...
ref                                                        34
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9193382 52.3892587, 4.9192...
Name: 9, dtype: object'
$ws.Range("A12").Value = 'osm_id                                                8456882
route                                                     bus
to                                  Amsterdam Olof Palmeplein
name        Bus 35: Amsterdam Molenwijk => Amsterdam Olof ...
ref                                                        35
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8936702 52.4180657, 4.8937...
Name: 10, dtype: object'
$ws.Range("A13").Value = 'osm_id                                                8456944
route                                                     bus
to                                        Amsterdam Molenwijk
name        Bus 35: Amsterdam Olof Palmeplein => Amsterdam...
ref                                                        35
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9401433 52.396837, 4.94000...
Name: 11, dtype: object'
$ws.Range("A18").Value = 'osm_id                                                8461307
route                                                     bus
to                                Amsterdam Metaalbewerkerweg
name        Bus 38: Amsterdam Station Noord => Amsterdam B...
ref                                                        38
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9333669 52.4019911, 4.9332...
Name: 16, dtype: object'
$ws.Range("A20").Value = 'osm_id                                                4539112
route                                                     bus
to                                    Amsterdam Amstelstation
name        Bus 40: Amsterdam Muiderpoortstation => Amster...
ref                                                        40
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9337209 52.3610843, 4.9339...
Name: 18, dtype: object'
$ws.Range("A21").Value = 'osm_id                                                  73789
route                                                     bus
to                               Amsterdam Muiderpoortstation
name        Bus 40: Amsterdam Amstelstation => Amsterdam M...
ref                                                        40
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9186931 52.3461675, 4.9186...
Name: 19, dtype: object'
$ws.Range("A22").Value = 'osm_id                                                 365785
route                                                     bus
to                             Amsterdam, Station Holendrecht
name        Bus 41: Amsterdam Muiderpoortstation => Amster...
ref                                                        41
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9771905 52.3236949, 4.9773...
Name: 20, dtype: object'
$ws.Range("A26").Value = 'osm_id                                                4542861
route                                                     bus
to                            Amsterdam Station Bijlmer ArenA
name        Bus 44: Diemen Noord => Amsterdam Station Bijl...
ref                                                        44
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9663154 52.3440833, 4.9663...
Name: 24, dtype: object'
$ws.Range("A28").Value = 'osm_id                                                 365789
route                                                     bus
to                            Amsterdam Station Bijlmer ArenA
name        Bus 47: Amsterdam Station Holendrecht => Amste...
ref                                                        47
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9492523 52.3131284, 4.9493...
Name: 26, dtype: object'
$ws.Range("A29").Value = 'osm_id                                                 365788
route                                                     bus
to                              Amsterdam Station Holendrecht
name        Bus 47: Amsterdam Station Bijlmer ArenA => Ams...
ref                                                        47
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.947369 52.3111419, 4.94752...
Name: 27, dtype: object'
$ws.Range("A31").Value = 'osm_id                                                4543993
route                                                     bus
to                                    Amsterdam, Koivistokade
name        Bus 48: Amsterdam Centraal Station => Amsterda...
ref                                                        48
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9006338 52.3795871, 4.9002...
Name: 29, dtype: object'
$ws.Range("A32").Value = 'osm_id                                                 365790
route                                                     bus
to                                              Weesp Station
name        Bus 49: Amsterdam Station Bijlmer ArenA => Wee...
ref                                                        49
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9474141 52.3110955, 4.9475...
Name: 30, dtype: object'
$ws.Range("A33").Value = 'osm_id                                                4547026
route                                                     bus
to                            Amsterdam Station Bijlmer ArenA
name        Bus 49: Weesp Station => Amsterdam Station Bij...
ref                                                        49
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.9919942 52.3142796, 4.9921...
Name: 31, dtype: object'
$ws.Range("A34").Value = 'osm_id                                                 375858
route                                                     bus
to                                Amsterdam Osdorpplein Noord
name        Bus 61: Amsterdam Station Sloterdijk => Amster...
ref                                                        61
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8387533 52.3881185, 4.8386...
Name: 32, dtype: object'
$ws.Range("A35").Value = 'osm_id                                                 365437
route                                                     bus
to                               Amsterdam Station Sloterdijk
name        Bus 61: Amsterdam Osdorpplein => Amsterdam Sta...
ref                                                        61
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8051399 52.3602896, 4.8047...
Name: 33, dtype: object'
$ws.Range("A38").Value = 'osm_id                                                4552870
route                                                     bus
to                                 Amsterdam Station Lelylaan
name        Bus 63: Osdorp De Aker => Amsterdam Station Le...
ref                                                        63
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.7745693 52.3543998, 4.7753...
Name: 36, dtype: object'
$ws.Range("A39").Value = 'osm_id                                                 375860
route                                                     bus
to                                             Osdorp De Aker
name        Bus 63: Amsterdam Station Lelylaan => Osdorp D...
ref                                                        63
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8343566 52.3566723, 4.8343...
Name: 37, dtype: object'
$ws.Range("A41").Value = 'osm_id                                                4552871
route                                                     bus
to                                       Amsterdam, KNSM-laan
name        Bus 65: Amsterdam Station Zuid => Amsterdam KN...
ref                                                        65
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8768122 52.3431691, 4.8768...
Name: 39, dtype: object'
$ws.Range("A44").Value = 'osm_id                                                5860193
route                                                     bus
to                             Amsterdam, John M. Keynesplein
name        Bus 68: Amsterdam Henk Sneevlietweg => Amsterd...
ref                                                        68
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8330587 52.3457424, 4.8328...
Name: 42, dtype: object'
$ws.Range("A45").Value = 'osm_id                                                6096663
route                                                     bus
to                               Amsterdam, Henk Sneevlietweg
name        Bus 68: Amsterdam Riekerpolder => Amsterdam He...
ref                                                        68
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.824444 52.340413, 4.824449...
Name: 43, dtype: object'
$ws.Range("A46").Value = 'osm_id                                               11720818
route                                                     bus
to                               Amsterdam, Henk Sneevlietweg
name        Bus 68: Amsterdam Riekerpolder => Amsterdam He...
ref                                                        68
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8252211 52.340147, 4.82543...
Name: 44, dtype: object'
$ws.Range("A47").Value = 'osm_id                                               11720819
route                                                     bus
to                             Amsterdam, John M. Keynesplein
name        Bus 68: Amsterdam Henk Sneevlietweg => Amsterd...
ref                                                        68
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8330587 52.3457424, 4.8328...
Name: 45, dtype: object'
$ws.Range("A49").Value = 'osm_id                                                 375862
route                                                     bus
to                                       Amsterdam, Abberdaan
name        Bus 231: Amsterdam Station Sloterdijk => Amste...
ref                                                       231
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8388353 52.3874811, 4.8387...
Name: 47, dtype: object'
$ws.Range("A56").Value = 'osm_id                                               11730321
route                                                     bus
to                             Amsterdam, John M. Keynesplein
name        Bus 267: Amsterdam Sloten => Amsterdam Riekerp...
ref                                                       267
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8143076 52.3354966, 4.8143...
Name: 54, dtype: object'
$ws.Range("A57").Value = 'osm_id                                                 375866
route                                                     bus
to                                    Schiphol, Airport/Plaza
name        Bus 369: Amsterdam Station Sloterdijk => Schip...
ref                                                       369
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.838741 52.38822, 4.838631 ...
Name: 55, dtype: object'
$ws.Range("A59").Value = 'osm_id                                                6900423
route                                                     bus
to                               Amsterdam Gustav Mahlerplein
name        Bus 461: Amsterdam Gelderlandplein NO => Amste...
ref                                                       461
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8791445 52.3318593, 4.8791...
Name: 57, dtype: object'
$ws.Range("A60").Value = 'osm_id                                                6900422
route                                                     bus
to                               Amsterdam Gelderlandplein NO
name        Bus 461: Amsterdam Gustav Mahlerplein => Amste...
ref                                                       461
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8731084 52.336651, 4.87323...
Name: 58, dtype: object'
$ws.Range("A61").Value = 'osm_id                                                6900421
route                                                     bus
to                                       Amsterdam, Bolestein
name        Bus 463: Amsterdam Gelderlandplein Oost => Ams...
ref                                                       463
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8794156 52.3308129, 4.8794...
Name: 59, dtype: object'
$ws.Range("A62").Value = 'osm_id                                                6900419
route                                                     bus
to                           Amsterdam, Gelderlandplein Noord
name        Bus 463: Amsterdam Bolestein => Amsterdam Geld...
ref                                                       463
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8596699 52.3249815, 4.8596...
Name: 60, dtype: object'
$ws.Range("A63").Value = 'osm_id                                               11730318
route                                                     bus
to                                    Amsterdam, De Boelelaan
name        Bus 464: Amsterdam De Boelelaan/De Klencke => ...
ref                                                       464
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8873203 52.3343831, 4.8875...
Name: 61, dtype: object'
$ws.Range("A64").Value = 'osm_id                                               11730319
route                                                     bus
to                                    Amsterdam, De Boelelaan
name        Bus 464: Amsterdam Gelderlandplein NO => Amste...
ref                                                       464
network                                Stadsvervoer Amsterdam
service                                                  None
geometry    MULTILINESTRING ((4.8791445 52.3318593, 4.8791...
Name: 62, dtype: object'
